$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 754.3333
$ws.Range("I6").Value = 805.2
$ws.Range("K6").Value = 2415.6
$ws.Range("M6").Value = -2303.6
$ws.Range("H12").Value = 214.21428
$ws.Range("I12").Value = 235
$ws.Range("J12").Value = 138
$ws.Range("K12").Value = 235
$ws.Range("L12").Value = 138
$ws.Range("M12").Value = -65
$ws.Range("N12").Value = -478
$ws.Range("H29").Value = 6126.75
$ws.Range("J29").Value = 7335.6665
$ws.Range("L29").Value = 22006.9995
$ws.Range("N29").Value = -22568.9995
$ws.Range("H33").Value = 453.45456
$ws.Range("I33").Value = 430
$ws.Range("K33").Value = 430
$ws.Range("M33").Value = -201
$ws.Range("H38").Value = 1091.4286
$ws.Range("I38").Value = 232.08333
$ws.Range("J38").Value = 6247.5
$ws.Range("K38").Value = 696.24999
$ws.Range("L38").Value = 18742.5
$ws.Range("M38").Value = -324.24999
$ws.Range("N38").Value = -19486.5
$ws.Range("H86").Value = 4124.4443
$ws.Range("I86").Value = 3611.875
$ws.Range("J86").Value = 4534.5
$ws.Range("K86").Value = 3611.875
$ws.Range("L86").Value = 4534.5
$ws.Range("M86").Value = -2488.875
$ws.Range("N86").Value = -6780.5
$ws.Range("H89").Value = 4124.4443
$ws.Range("I89").Value = 3611.875
$ws.Range("J89").Value = 4534.5
$ws.Range("K89").Value = 18059.375
$ws.Range("L89").Value = 22672.5
$ws.Range("M89").Value = -12443.375
$ws.Range("N89").Value = -33904.5
$ws.Range("H127").Value = 1481
$ws.Range("I127").Value = 568.3333
$ws.Range("J127").Value = 2850
$ws.Range("K127").Value = 1704.9999
$ws.Range("L127").Value = 8550
$ws.Range("M127").Value = 3255.0001
$ws.Range("N127").Value = -18470
$ws.Range("H141").Value = 2591.0557
$ws.Range("I141").Value = 2027.8572
$ws.Range("K141").Value = 6083.571599999999
$ws.Range("M141").Value = -903.5715999999993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7463.6665
$ws.Range("I2").Value = 685.6
$ws.Range("K2").Value = 685.6
$ws.Range("M2").Value = -572.6
$ws.Range("H30").Value = 1300
$ws.Range("I30").Value = 1500
$ws.Range("J30").Value = 1100
$ws.Range("K30").Value = 1500
$ws.Range("L30").Value = 1100
$ws.Range("M30").Value = -1350
$ws.Range("N30").Value = -1400
$ws.Range("H74").Value = 2768.6785
$ws.Range("I74").Value = 1906.625
$ws.Range("K74").Value = 1906.625
$ws.Range("M74").Value = -1032.625
$ws.Range("H77").Value = 2768.6785
$ws.Range("I77").Value = 1906.625
$ws.Range("K77").Value = 9533.125
$ws.Range("M77").Value = -5165.125
$ws.Range("H116").Value = 7463.6665
$ws.Range("I116").Value = 685.6
$ws.Range("K116").Value = 685.6
$ws.Range("M116").Value = 1608.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7463.6665
$ws.Range("I3").Value = 685.6
$ws.Range("K3").Value = 685.6
$ws.Range("M3").Value = -571.6
$ws.Range("H86").Value = 5831.2
$ws.Range("I86").Value = 5831.2
$ws.Range("K86").Value = 5831.2
$ws.Range("M86").Value = -4708.2
$ws.Range("H89").Value = 5831.2
$ws.Range("I89").Value = 5831.2
$ws.Range("K89").Value = 29156
$ws.Range("M89").Value = -23540
$ws.Range("H105").Value = 13788.741
$ws.Range("I105").Value = 11916.55
$ws.Range("J105").Value = 19137.857
$ws.Range("K105").Value = 11916.55
$ws.Range("L105").Value = 19137.857
$ws.Range("M105").Value = -10169.55
$ws.Range("N105").Value = -22631.857
$ws.Range("H134").Value = 1999.3636
$ws.Range("I134").Value = 1573.95
$ws.Range("K134").Value = 4721.85
$ws.Range("M134").Value = -2186.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 431.8
$ws.Range("J7").Value = 472
$ws.Range("L7").Value = 472
$ws.Range("N7").Value = -698
$ws.Range("H86").Value = 3826.3333
$ws.Range("I86").Value = 4099.5386
$ws.Range("J86").Value = 3116
$ws.Range("K86").Value = 4099.5386
$ws.Range("L86").Value = 3116
$ws.Range("M86").Value = -2976.5386
$ws.Range("N86").Value = -5362
$ws.Range("H89").Value = 3826.3333
$ws.Range("I89").Value = 4099.5386
$ws.Range("J89").Value = 3116
$ws.Range("K89").Value = 20497.693
$ws.Range("L89").Value = 15580
$ws.Range("M89").Value = -14881.693
$ws.Range("N89").Value = -26812

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 7074852
$ws.Range("I34").Value = 12965562
$ws.Range("J34").Value = 5999.8
$ws.Range("K34").Value = 38896686
$ws.Range("L34").Value = 17999.4
$ws.Range("M34").Value = -38896602
$ws.Range("N34").Value = -18167.4
$ws.Range("H39").Value = 1798.3334
$ws.Range("J39").Value = 2900
$ws.Range("L39").Value = 8700
$ws.Range("N39").Value = -9288
$ws.Range("H55").Value = 1852.8462
$ws.Range("I55").Value = 1590.5834
$ws.Range("J55").Value = 5000
$ws.Range("K55").Value = 4771.7502
$ws.Range("L55").Value = 15000
$ws.Range("M55").Value = -4594.7502
$ws.Range("N55").Value = -15354
$ws.Range("H132").Value = 4957.143
$ws.Range("I132").Value = 4156.077
$ws.Range("J132").Value = 5651.4
$ws.Range("K132").Value = 37404.693
$ws.Range("L132").Value = 50862.6
$ws.Range("M132").Value = -34874.693
$ws.Range("N132").Value = -55922.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 339893.2
$ws.Range("I80").Value = 557459.9
$ws.Range("J80").Value = 13543.167
$ws.Range("K80").Value = 557459.9
$ws.Range("L80").Value = 13543.167
$ws.Range("M80").Value = -556461.9
$ws.Range("N80").Value = -15539.167
$ws.Range("H83").Value = 339893.2
$ws.Range("I83").Value = 557459.9
$ws.Range("J83").Value = 13543.167
$ws.Range("K83").Value = 2787299.5
$ws.Range("L83").Value = 67715.83499999999
$ws.Range("M83").Value = -2782307.5
$ws.Range("N83").Value = -77699.83499999999
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820
$ws.Range("H132").Value = 2380.7847
$ws.Range("I132").Value = 2111.8525
$ws.Range("K132").Value = 6335.5575
$ws.Range("M132").Value = -3805.5575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 100639.75
$ws.Range("J128").Value = 100639.75
$ws.Range("L128").Value = 100639.75
$ws.Range("N128").Value = -110599.75
